{"js": "const replacements = [\n  [\"473\u00d73=1419\", \"383\u00d72=766\"],\n  [\"781\u00d73=2343\", \"530\u00d78=4240\"],\n  [\"768\u00d77=5376\", \"705\u00d72=1410\"],\n  [\"999\u00d75=4995\", \"425\u00d78=3400\"],\n  [\"170\u00d76=1020\", \"146\u00d76=876\"],\n  [\"441\u00d78=3528\", \"296\u00d73=888\"],\n  [\"822\u00d79=7398\", \"333\u00d79=2997\"],\n  [\"363\u00d77=2541\", \"382\u00d75=1910\"],\n  [\"517\u00d79=4653\", \"607\u00d77=4249\"],\n  [\"561\u00d73=1683\", \"661\u00d76=3966\"],\n  [\"313\u00d78=2504\", \"883\u00d72=1766\"],\n  [\"134\u00d72=268\", \"424\u00d78=3392\"],\n  [\"781\u00d77=5467\", \"341\u00d72=682\"],\n  [\"680\u00d78=5440\", \"326\u00d78=2608\"],\n  [\"829\u00d74=3316\", \"259\u00d79=2331\"],\n  [\"225\u00d75=1125\", \"448\u00d79=4032\"],\n  [\"857\u00d79=7713\", \"308\u00d77=2156\"],\n  [\"798\u00d79=7182\", \"503\u00d79=4527\"],\n  [\"271\u00d75=1355\", \"133\u00d77=931\"],\n  [\"671\u00d75=3355\", \"782\u00d78=6256\"],\n  [\"268\u00d78=2144\", \"925\u00d74=3700\"],\n  [\"653\u00d74=2612\", \"220\u00d74=880\"],\n  [\"987\u00d77=6909\", \"308\u00d79=2772\"],\n  [\"929\u00d78=7432\", \"292\u00d79=2628\"],\n  [\"323\u00d78=2584\", \"118\u00d74=472\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"473\u00d73=1419\"\n$find.Replacement.Text = \"383\u00d72=766\"\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"781\u00d73=2343\"\n$find.Replacement.Text = \"530\u00d78=4240\"\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"768\u00d77=5376\"\n$find.Replacement.Text = \"705\u00d72=1410\"\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"999\u00d75=4995\"\n$find.Replacement.Text = \"425\u00d78=3400\"\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"170\u00d76=1020\"\n$find.Replacement.Text = \"146\u00d76=876\"\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"441\u00d78=3528\"\n$find.Replacement.Text = \"296\u00d73=888\"\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"822\u00d79=7398\"\n$find.Replacement.Text = \"333\u00d79=2997\"\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"363\u00d77=2541\"\n$find.Replacement.Text = \"382\u00d75=1910\"\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"517\u00d79=4653\"\n$find.Replacement.Text = \"607\u00d77=4249\"\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"561\u00d73=1683\"\n$find.Replacement.Text = \"661\u00d76=3966\"\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"313\u00d78=2504\"\n$find.Replacement.Text = \"883\u00d72=1766\"\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"134\u00d72=268\"\n$find.Replacement.Text = \"424\u00d78=3392\"\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"781\u00d77=5467\"\n$find.Replacement.Text = \"341\u00d72=682\"\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"680\u00d78=5440\"\n$find.Replacement.Text = \"326\u00d78=2608\"\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"829\u00d74=3316\"\n$find.Replacement.Text = \"259\u00d79=2331\"\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"225\u00d75=1125\"\n$find.Replacement.Text = \"448\u00d79=4032\"\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"857\u00d79=7713\"\n$find.Replacement.Text = \"308\u00d77=2156\"\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"798\u00d79=7182\"\n$find.Replacement.Text = \"503\u00d79=4527\"\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"271\u00d75=1355\"\n$find.Replacement.Text = \"133\u00d77=931\"\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"671\u00d75=3355\"\n$find.Replacement.Text = \"782\u00d78=6256\"\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"268\u00d78=2144\"\n$find.Replacement.Text = \"925\u00d74=3700\"\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"653\u00d74=2612\"\n$find.Replacement.Text = \"220\u00d74=880\"\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"987\u00d77=6909\"\n$find.Replacement.Text = \"308\u00d79=2772\"\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"929\u00d78=7432\"\n$find.Replacement.Text = \"292\u00d79=2628\"\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"323\u00d78=2584\"\n$find.Replacement.Text = \"118\u00d74=472\"\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n"}
